# Saldo.xlsx update:
#  - Remove the "Gustavo" row (004565146 / Gustavo / 6598.22)
#  - Replace the "Igor" row (008054285 / Igor / 5007.95), which shifts up to
#    take Gustavo's old position, with new data (004368628 / Camila / 2116.76)
#  - Remove the old "Camila" row (004368628 / Camila / 179.88) further down,
#    which shifts up by one row once the Gustavo row above it is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the entire "Gustavo" row (row 5). Everything below shifts up by one.
$ws.Rows.Item(5).Delete()

# 2) The old Igor row (was row 7) is now row 6 after the shift above.
#    Overwrite it in place with the new Camila data. The account number must
#    stay text (keeping its leading zeros), so the cell is formatted as
#    Text before the value is assigned.
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "004368628"
$ws.Cells.Item(6, 2).Value = "Camila"
$ws.Cells.Item(6, 3).Value = 2116.76

# 3) The old Camila/179.88 row (was row 48) is now row 47 after the shift
#    above. Delete it entirely.
$ws.Rows.Item(47).Delete()
